$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 59, pushing existing rows 59-63 down to 61-65
$ws.Rows.Item(59).Insert()
$ws.Rows.Item(59).Insert()

# Row 59: new weekly data (Primera quality)
$ws.Cells.Item(59, 1).Value = 9
$ws.Cells.Item(59, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(59, 3).Value = "Metropolitana"
$ws.Cells.Item(59, 4).Value = 44578
$ws.Cells.Item(59, 5).Value = 13
$ws.Cells.Item(59, 6).Value = 100114002
$ws.Cells.Item(59, 7).Value = "Camote"
$ws.Cells.Item(59, 8).Value = "Sin especificar"
$ws.Cells.Item(59, 9).Value = "Primera"
$ws.Cells.Item(59, 10).Value = 1240
$ws.Cells.Item(59, 11).Value = 10000
$ws.Cells.Item(59, 12).Value = 11000
$ws.Cells.Item(59, 13).Value = 10500
$ws.Cells.Item(59, 14).Value = "`$/malla 18 kilos"
$ws.Cells.Item(59, 15).Value = "Perú"
$ws.Cells.Item(59, 16).Value = 583
$ws.Cells.Item(59, 17).Value = 18
$ws.Cells.Item(59, 18).Value = "Hortaliza"

# Row 60: new weekly data (Segunda quality)
$ws.Cells.Item(60, 1).Value = 9
$ws.Cells.Item(60, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(60, 3).Value = "Metropolitana"
$ws.Cells.Item(60, 4).Value = 44578
$ws.Cells.Item(60, 5).Value = 13
$ws.Cells.Item(60, 6).Value = 100114002
$ws.Cells.Item(60, 7).Value = "Camote"
$ws.Cells.Item(60, 8).Value = "Sin especificar"
$ws.Cells.Item(60, 9).Value = "Segunda"
$ws.Cells.Item(60, 10).Value = 610
$ws.Cells.Item(60, 11).Value = 9000
$ws.Cells.Item(60, 12).Value = 9000
$ws.Cells.Item(60, 13).Value = 9000
$ws.Cells.Item(60, 14).Value = "`$/malla 18 kilos"
$ws.Cells.Item(60, 15).Value = "Perú"
$ws.Cells.Item(60, 16).Value = 500
$ws.Cells.Item(60, 17).Value = 18
$ws.Cells.Item(60, 18).Value = "Hortaliza"
